# Update team-specific transition-probability matrix values on Sheet1
# (Colorado_B team simulation matrix) with refreshed probabilities.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1826697892271663
$ws.Range("C2").Value = 0.5667447306791569
$ws.Range("J2").Value = 0.00936768149882904
$ws.Range("P2").Value = 0.1334894613583138
$ws.Range("S2").Value = 0.107728337236534
$ws.Range("B3").Value = 0.007936507936507936
$ws.Range("C3").Value = 0.03571428571428571
$ws.Range("J3").Value = 0.02777777777777778
$ws.Range("P3").Value = 0.7579365079365079
$ws.Range("S3").Value = 0.1706349206349206
$ws.Range("J4").Value = 0.07843137254901961
$ws.Range("P4").Value = 0.7254901960784313
$ws.Range("S4").Value = 0.196078431372549
$ws.Range("P5").Value = 0.8
$ws.Range("S5").Value = 0.2
$ws.Range("B6").Value = 0.06666666666666667
$ws.Range("D6").Value = 0.005797101449275362
$ws.Range("F6").Value = 0.04057971014492753
$ws.Range("J6").Value = 0.2956521739130435
$ws.Range("O6").Value = 0.01449275362318841
$ws.Range("Q6").Value = 0.1884057971014493
$ws.Range("R6").Value = 0.05797101449275362
$ws.Range("S6").Value = 0.3304347826086956
$ws.Range("B7").Value = 0.120253164556962
$ws.Range("D7").Value = 0.02848101265822785
$ws.Range("F7").Value = 0.05696202531645569
$ws.Range("J7").Value = 0.1360759493670886
$ws.Range("O7").Value = 0.02531645569620253
$ws.Range("Q7").Value = 0.1772151898734177
$ws.Range("R7").Value = 0.06012658227848101
$ws.Range("S7").Value = 0.3955696202531646
$ws.Range("B8").Value = 0.08781127129750983
$ws.Range("D8").Value = 0.01310615989515072
$ws.Range("E8").Value = 0.002621231979030144
$ws.Range("F8").Value = 0.07863695937090433
$ws.Range("J8").Value = 0.1022280471821756
$ws.Range("O8").Value = 0.01703800786369594
$ws.Range("Q8").Value = 0.1939711664482307
$ws.Range("R8").Value = 0.07077326343381389
$ws.Range("S8").Value = 0.4338138925294889
$ws.Range("B9").Value = 0.09935897435897435
$ws.Range("D9").Value = 0.009615384615384616
$ws.Range("F9").Value = 0.04166666666666666
$ws.Range("J9").Value = 0.125
$ws.Range("O9").Value = 0.02564102564102564
$ws.Range("Q9").Value = 0.1826923076923077
$ws.Range("R9").Value = 0.0608974358974359
$ws.Range("S9").Value = 0.4551282051282051
$ws.Range("B10").Value = 0.1052631578947368
$ws.Range("D10").Value = 0.01659038901601831
$ws.Range("E10").Value = 0.002288329519450801
$ws.Range("F10").Value = 0.07608695652173914
$ws.Range("J10").Value = 0.108695652173913
$ws.Range("O10").Value = 0.0137299771167048
$ws.Range("Q10").Value = 0.2053775743707094
$ws.Range("R10").Value = 0.07608695652173914
$ws.Range("S10").Value = 0.3958810068649886
$ws.Range("G11").Value = 0.1331658291457286
$ws.Range("J11").Value = 0.06281407035175879
$ws.Range("K11").Value = 0.1507537688442211
$ws.Range("L11").Value = 0.6407035175879398
$ws.Range("S11").Value = 0.01256281407035176
$ws.Range("G12").Value = 0.7667844522968198
$ws.Range("J12").Value = 0.1696113074204947
$ws.Range("K12").Value = 0.007067137809187279
$ws.Range("L12").Value = 0.02473498233215548
$ws.Range("S12").Value = 0.03180212014134275
$ws.Range("G13").Value = 0.7108433734939759
$ws.Range("J13").Value = 0.2771084337349398
$ws.Range("S13").Value = 0.01204819277108434
$ws.Range("F15").Value = 0.02064896755162242
$ws.Range("H15").Value = 0.1799410029498525
$ws.Range("I15").Value = 0.08259587020648967
$ws.Range("J15").Value = 0.3185840707964602
$ws.Range("K15").Value = 0.07964601769911504
$ws.Range("M15").Value = 0.02359882005899705
$ws.Range("O15").Value = 0.05309734513274336
$ws.Range("S15").Value = 0.2418879056047198
$ws.Range("F16").Value = 0.0176056338028169
$ws.Range("H16").Value = 0.2253521126760563
$ws.Range("I16").Value = 0.09507042253521127
$ws.Range("J16").Value = 0.3345070422535211
$ws.Range("K16").Value = 0.1056338028169014
$ws.Range("M16").Value = 0.04577464788732395
$ws.Range("O16").Value = 0.08450704225352113
$ws.Range("S16").Value = 0.09154929577464789
$ws.Range("F17").Value = 0.0218978102189781
$ws.Range("H17").Value = 0.2102189781021898
$ws.Range("I17").Value = 0.08759124087591241
$ws.Range("J17").Value = 0.3839416058394161
$ws.Range("K17").Value = 0.09781021897810219
$ws.Range("M17").Value = 0.01605839416058394
$ws.Range("O17").Value = 0.07153284671532846
$ws.Range("S17").Value = 0.1109489051094891
$ws.Range("F18").Value = 0.01239669421487603
$ws.Range("H18").Value = 0.2355371900826446
$ws.Range("I18").Value = 0.05371900826446281
$ws.Range("J18").Value = 0.371900826446281
$ws.Range("K18").Value = 0.1446280991735537
$ws.Range("M18").Value = 0.02892561983471074
$ws.Range("N18").Value = 0.004132231404958678
$ws.Range("O18").Value = 0.07024793388429752
$ws.Range("S18").Value = 0.07851239669421488
$ws.Range("F19").Value = 0.02218700475435816
$ws.Range("H19").Value = 0.2282091917591125
$ws.Range("I19").Value = 0.09561542525092445
$ws.Range("J19").Value = 0.358161648177496
$ws.Range("K19").Value = 0.09561542525092445
$ws.Range("M19").Value = 0.02377179080824089
$ws.Range("O19").Value = 0.07237189646064449
$ws.Range("S19").Value = 0.109878499735869
